$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.636.46'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.847.81'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.40'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4279'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3630'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.84'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07314'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.17%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.65'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.863.04'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.324'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.516'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06904'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '79.87'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000009034'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.94%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.34'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '27.664.33'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.959'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.38'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.093.42'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.990'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.20'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.70'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '121.50'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +8.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.289'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.25%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08905'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7619'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.978'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.25%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.63%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05405'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.71%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01931'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.815'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5073'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1655'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.770'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.352'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.76%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.06547'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.31'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '105.13'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4680'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.76%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.618'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.44'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.35%  '
